# Apply updated "dSF" (column F) values after a data repull / mean recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = 4
    4  = 3
    5  = 3
    6  = -2
    7  = -1
    8  = -1
    9  = 2
    10 = 4
    12 = 3
    13 = 6
    14 = -5
    15 = -1
    16 = 5
    17 = 2
    18 = 3
    19 = -3
    21 = 9
    23 = -5
    24 = 2
    25 = -6
    26 = 2
    27 = 1
    29 = 4
    30 = -1
    31 = -6
    32 = 3
    35 = -6
    36 = -1
    37 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
